# Insert Stefan Lyocsa's author line and publications 18-23 immediately
# after publication #17, before the "Lennart John Baals" paragraph.

$d = $word.ActiveDocument

$target = $d.Content.Duplicate
$found = $target.Find.Execute(
    "17. `"Predicting Retail Customers' Distress: Early Warning Systems and Machine Learning Applications`". SSRN. DOI: 10.2139/ssrn.4730470",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find publication #17 text"
}

# Collapse the found range to its end point (wdCollapseEnd = 0), right
# after "...4730470".
$target.Collapse(0)

$lines = @(
    "Stefan Lyocsa (ORCID: 0000-0002-8380-181X):",
    "18. `"Macroeconomic environment and the future performance of loans: Evidence from three peer-to-peer platforms`". International Review of Financial Analysis. DOI: 10.1016/j.irfa.2024.103416",
    "19. `"What drives the uranium sector risk? The role of attention, economic and geopolitical uncertainty`". Energy Economics. DOI: 10.1016/j.eneco.2024.107980",
    "20. `"Forecasting of clean energy market volatility: The role of oil and the technology sector`". Energy Economics. DOI: 10.1016/j.eneco.2024.107451",
    "21. `"A Fuzzy Framework for Realized Volatility Prediction`" (2025). SSRN.",
    "22. `"Alpha-threshold networks in credit risk models`" (2025). SSRN.",
    "23. `"Do hurricanes cause storm on the stock market?`" (2025). SSRN."
)

# Each new line becomes its own paragraph: prefix every line (including
# the first) with a paragraph mark so they land as separate <w:p> blocks
# after the existing "17. ..." paragraph.
$insertText = "`r" + ($lines -join "`r")
$target.InsertAfter($insertText)
